$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 (I0) and J1 (IF), matching the formatting of the
# existing header cells (bold font, border, centered alignment).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Row data for the new I (I0) and J (IF) columns: row index, I value, J value
$data = @(
    @(2, 7, 7),
    @(3, 8, 8),
    @(4, 8, 8),
    @(5, 8, 8),
    @(6, 8, 8),
    @(7, 9, 9),
    @(8, 6, 6),
    @(9, 8, 8),
    @(10, 6, 6),
    @(11, 7, 7),
    @(12, 7, 8),
    @(13, 8, 8),
    @(14, 7, 7),
    @(15, 8, 8),
    @(16, 7, 8),
    @(17, 8, 9),
    @(18, 8, 8),
    @(19, 7, 8),
    @(20, 8, 8),
    @(21, 8, 8),
    @(22, 8, 9),
    @(23, 8, 8),
    @(24, 8, 8),
    @(25, 8, 8),
    @(26, 7, 7),
    @(27, 8, 8),
    @(28, 8, 8),
    @(29, 7, 7),
    @(30, 7, 7),
    @(31, 7, 7),
    @(32, 9, 9),
    @(33, 8, 8),
    @(34, 9, 9),
    @(35, 7, 7),
    @(36, 7, 7),
    @(37, 7, 7),
    @(38, 7, 7),
    @(39, 10, 10),
    @(40, 6, 6),
    @(41, 9, 9),
    @(42, 7, 7),
    @(43, 7, 7),
    @(44, 7, 7),
    @(45, 7, 8),
    @(46, 6, 7),
    @(47, 7, 8),
    @(48, 8, 8),
    @(49, 5, 5),
    @(50, 9, 9),
    @(51, 8, 8),
    @(52, 5, 6),
    @(53, 7, 7),
    @(54, 7, 7),
    @(55, 6, 7),
    @(56, 6, 7),
    @(57, 7, 7),
    @(58, 7, 7),
    @(59, 7, 7),
    @(60, 8, 8),
    @(61, 6, 6),
    @(62, 8, 8),
    @(63, 8, 8),
    @(64, 8, 8),
    @(65, 7, 7),
    @(66, 8, 8),
    @(67, 7, 7),
    @(68, 9, 9),
    @(69, 6, 6),
    @(70, 7, 7),
    @(71, 6, 7),
    @(72, 6, 7),
    @(73, 7, 7),
    @(74, 8, 9),
    @(75, 8, 8),
    @(76, 6, 6),
    @(77, 6, 6),
    @(78, 8, 8),
    @(79, 8, 8),
    @(80, 6, 6),
    @(81, 8, 8),
    @(82, 7, 7),
    @(83, 6, 6),
    @(84, 6, 6),
    @(85, 8, 8),
    @(86, 6, 6),
    @(87, 5, 5),
    @(88, 6, 6),
    @(89, 6, 7),
    @(90, 5, 5),
    @(91, 5, 5),
    @(92, 5, 5),
    @(93, 6, 6),
    @(94, 6, 6)
)

foreach ($item in $data) {
    $row = $item[0]
    $iVal = $item[1]
    $jVal = $item[2]
    $ws.Cells.Item($row, 9).Value = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
}
